$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe
# (Excel will not reinterpret them as numbers).
$ws.Range("D2").Value = '67.655.50'
$ws.Range("E2").Value = '  +0.34%  '

$ws.Range("D3").Value = '3.508.93'
$ws.Range("E3").Value = '  -0.42%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("E5").Value = '  -0.89%  '

$ws.Range("E6").Value = '  +0.46%  '

$ws.Range("D7").Value = '3.506.09'
$ws.Range("E7").Value = '  -0.50%  '

$ws.Range("E9").Value = '  +1.08%  '

$ws.Range("E10").Value = '  +2.70%  '

$ws.Range("E11").Value = '  +7.69%  '

$ws.Range("E12").Value = '  +1.62%  '

$ws.Range("E13").Value = '  -1.89%  '

$ws.Range("E14").Value = '  +0.31%  '

$ws.Range("D15").Value = '4.101.84'
$ws.Range("E15").Value = '  -0.47%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.511.76'
$ws.Range("E16").Value = '  -0.51%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '67.571.61'
$ws.Range("E17").Value = '  +0.18%  '

$ws.Range("E18").Value = '  -0.63%  '

$ws.Range("E19").Value = '  +1.81%  '

$ws.Range("E20").Value = '  +0.61%  '

$ws.Range("E21").Value = '  +3.34%  '

$ws.Range("E22").Value = '  +0.17%  '

$ws.Range("E23").Value = '  +0.80%  '

$ws.Range("E24").Value = '  +1.12%  '

$ws.Range("D25").Value = '3.648.04'
$ws.Range("E25").Value = '  -0.49%  '

$ws.Range("E26").Value = '  -3.62%  '

$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("E28").Value = '  +3.79%  '

$ws.Range("E29").Value = '  -1.74%  '

$ws.Range("E30").Value = '  +0.39%  '

$ws.Range("E31").Value = '  +5.42%  '

$ws.Range("E32").Value = '  +5.23%  '

$ws.Range("E33").Value = '  +0.14%  '

$ws.Range("E34").Value = '  -0.90%  '

$ws.Range("E35").Value = '  -0.58%  '

$ws.Range("E36").Value = '  +0.75%  '

$ws.Range("D37").Value = '3.498.38'
$ws.Range("E37").Value = '  -0.54%  '

$ws.Range("E38").Value = '  -0.33%  '

$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("E40").Value = '  +7.13%  '

$ws.Range("E41").Value = '  +1.11%  '

$ws.Range("E42").Value = '  -0.02%  '

$ws.Range("E43").Value = '  +1.75%  '

$ws.Range("E44").Value = '  +0.27%  '

$ws.Range("E45").Value = '  +1.33%  '

$ws.Range("E46").Value = '  +6.50%  '

$ws.Range("E47").Value = '  +2.79%  '

$ws.Range("E48").Value = '  +4.53%  '

$ws.Range("E49").Value = '  -3.20%  '

$ws.Range("E50").Value = '  +0.42%  '

$ws.Range("E51").Value = '  +1.94%  '

# Numeric-looking "Price" values (e.g. "1.00", "7.67") must stay as literal
# text, matching the original inlineStr cells. Setting .Value directly would
# make Excel auto-convert them to real numbers (and round-trip/precision
# would drift, e.g. "609.55" -> 609.54999999999995). We instead stage the
# text on a helper cell formatted as Text ("@"), copy it, and paste-special
# VALUES ONLY (xlPasteValues = -4163) into the target cell, which carries
# over the literal text without bringing the Text number format along -
# leaving the target cell style untouched, exactly like the diff expects.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$helper.Value = '1.00'
$helper.Copy()
$ws.Range("D4").PasteSpecial(-4163)

$helper.Value = '609.55'
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)

$helper.Value = '152.16'
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)

$helper.Value = '7.67'
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)

$helper.Value = '0.432'
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)

$helper.Value = '0.0000217'
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)

$helper.Value = '32.18'
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)

$helper.Value = '6.52'
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)

$helper.Value = '15.49'
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)

$helper.Value = '9.86'
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)

$helper.Value = '447.18'
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)

$helper.Value = '0.629'
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)

$helper.Value = '78.40'
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)

$helper.Value = '0.0000127'
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)

$helper.Value = '8.79'
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)

$helper.Value = '10.04'
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)

$helper.Value = '2.53'
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)

$helper.Value = '0.172'
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)

$helper.Value = '25.63'
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)

$helper.Value = '6.15'
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)

$helper.Value = '1.86'
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)

$helper.Value = '2.31'
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)

$helper.Value = '179.27'
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)

$helper.Value = '1.00'
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)

$helper.Value = '0.0898'
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)

$helper.Value = '5.45'
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)

$helper.Value = '0.894'
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)

$helper.Value = '30.22'
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)

$helper.Value = '46.46'
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)

$helper.Value = '2.56'
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)

$helper.Value = '0.253'
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)

# Clean up the helper cell so it leaves no trace in the saved workbook.
$helper.Clear()
